$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (e.g. "6.09", "560.42") are not auto-coerced into numbers by Excel,
# matching the original inline-string cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.069.55'
$ws.Range("E2").Value = '  +3.06%  '
$ws.Range("D3").Value = '3.062.57'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '560.42'
$ws.Range("E5").Value = '  +3.42%  '
$ws.Range("D6").Value = '144.11'
$ws.Range("E6").Value = '  +3.04%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.063.14'
$ws.Range("E8").Value = '  +2.24%  '
$ws.Range("D9").Value = '0.511'
$ws.Range("E9").Value = '  +4.78%  '
$ws.Range("E10").Value = '  +5.80%  '
$ws.Range("D11").Value = '6.09'
$ws.Range("E11").Value = '  -10.22%  '
$ws.Range("D12").Value = '0.484'
$ws.Range("E12").Value = '  +8.82%  '
$ws.Range("E13").Value = '  +5.60%  '
$ws.Range("D14").Value = '35.55'
$ws.Range("E14").Value = '  +4.71%  '
$ws.Range("D15").Value = '3.563.15'
$ws.Range("E15").Value = '  +2.29%  '
$ws.Range("D16").Value = '64.107.29'
$ws.Range("E16").Value = '  +2.93%  '
$ws.Range("D17").Value = '3.069.25'
$ws.Range("E17").Value = '  +2.04%  '
$ws.Range("E18").Value = '  +2.32%  '
$ws.Range("D19").Value = '6.77'
$ws.Range("E19").Value = '  +3.45%  '
$ws.Range("D20").Value = '478.56'
$ws.Range("E20").Value = '  +2.60%  '
$ws.Range("D21").Value = '13.99'
$ws.Range("E21").Value = '  +4.61%  '
$ws.Range("E22").Value = '  +4.62%  '
$ws.Range("D23").Value = '7.58'
$ws.Range("E23").Value = '  +5.62%  '
$ws.Range("D24").Value = '14.34'
$ws.Range("E24").Value = '  +14.06%  '
$ws.Range("D25").Value = '82.09'
$ws.Range("E25").Value = '  +3.43%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '2.80'
$ws.Range("E27").Value = '  +3.46%  '
$ws.Range("D28").Value = '8.02'
$ws.Range("E28").Value = '  +5.59%  '
$ws.Range("D29").Value = '2.04'
$ws.Range("E29").Value = '  +2.08%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").Value = '26.38'
$ws.Range("E31").Value = '  +3.81%  '
$ws.Range("E32").Value = '  +1.35%  '
$ws.Range("E33").Value = '  +4.72%  '
$ws.Range("D34").Value = '5.72'
$ws.Range("E34").Value = '  +3.22%  '
$ws.Range("E35").Value = '  +7.58%  '
$ws.Range("D36").Value = '54.99'
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("D37").Value = '0.0409'
$ws.Range("E37").Value = '  +5.18%  '
$ws.Range("D38").Value = '445.75'
$ws.Range("E38").Value = '  -0.76%  '
$ws.Range("D39").Value = '0.0812'
$ws.Range("E39").Value = '  +0.66%  '
$ws.Range("D40").Value = '2.85'
$ws.Range("E40").Value = '  +10.92%  '
$ws.Range("D41").Value = '2.995.90'
$ws.Range("E41").Value = '  +1.74%  '
$ws.Range("D42").Value = '8.25'
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("D43").Value = '0.115'
$ws.Range("E43").Value = '  +1.25%  '
$ws.Range("D44").Value = '27.99'
$ws.Range("E44").Value = '  +4.77%  '
$ws.Range("E45").Value = '  +6.34%  '
$ws.Range("D46").Value = '2.17'
$ws.Range("E46").Value = '  +8.40%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("E48").Value = '  +4.50%  '
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '0.0₃0520'
$ws.Range("E49").Value = '  +5.20%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '118.73'
$ws.Range("E50").Value = '  +3.10%  '
$ws.Range("E51").Value = '  +3.86%  '
